$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 520.9706
$ws.Range("I15").Value = 520.9706
$ws.Range("K15").Value = 1562.9118
$ws.Range("M15").Value = -1393.9118

$ws.Range("H33").Value = 680.8570999999999
$ws.Range("I33").Value = 544.4167
$ws.Range("J33").Value = 1499.5
$ws.Range("K33").Value = 544.4167
$ws.Range("L33").Value = 1499.5
$ws.Range("M33").Value = -315.4167
$ws.Range("N33").Value = -1957.5

$ws.Range("H51").Value = 3875

$ws.Range("H55").Value = 159.52632
$ws.Range("J55").Value = 263.27274
$ws.Range("L55").Value = 263.27274
$ws.Range("N55").Value = -691.27274

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""

$ws.Range("H106").Value = 797235.9
$ws.Range("I106").Value = 1390138
$ws.Range("K106").Value = 1390138
$ws.Range("M106").Value = -1389507

$ws.Range("H113").Value = 5323.4

$ws.Range("H116").Value = 1048807.2
$ws.Range("I116").Value = 7039.409
$ws.Range("K116").Value = 7039.409
$ws.Range("M116").Value = -3597.409

$ws.Range("H132").Value = 1726.8572
$ws.Range("I132").Value = 1304.6129
$ws.Range("K132").Value = 3913.8387
$ws.Range("M132").Value = -1383.8387

$ws.Range("H137").Value = 661376.1
$ws.Range("I137").Value = 2227.818
$ws.Range("J137").Value = 1320524.5
$ws.Range("K137").Value = 6683.454000000001
$ws.Range("L137").Value = 3961573.5
$ws.Range("M137").Value = -4133.454000000001
$ws.Range("N137").Value = -3966673.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""

$ws.Range("H61").Value = 1891.7391
$ws.Range("I61").Value = 1620.5
$ws.Range("K61").Value = 1620.5
$ws.Range("M61").Value = -1408.5

$ws.Range("H63").Value = 2459.7693
$ws.Range("I63").Value = 2543.3635
$ws.Range("K63").Value = 2543.3635
$ws.Range("M63").Value = -1857.3635

$ws.Range("H66").Value = 2459.7693
$ws.Range("I66").Value = 2543.3635
$ws.Range("K66").Value = 12716.8175
$ws.Range("M66").Value = -9284.817499999999

$ws.Range("H74").Value = 2191.75
$ws.Range("I74").Value = 1500
$ws.Range("J74").Value = 2422.3333
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 2422.3333
$ws.Range("M74").Value = -626
$ws.Range("N74").Value = -4170.3333

$ws.Range("H77").Value = 2191.75
$ws.Range("I77").Value = 1500
$ws.Range("J77").Value = 2422.3333
$ws.Range("K77").Value = 7500
$ws.Range("L77").Value = 12111.6665
$ws.Range("M77").Value = -3132
$ws.Range("N77").Value = -20847.6665

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""

$ws.Range("H102").Value = 78291.92999999999
$ws.Range("I102").Value = 85093.164
$ws.Range("K102").Value = 85093.164
$ws.Range("M102").Value = -83471.164

$ws.Range("H122").Value = 5647.8
$ws.Range("I122").Value = 6013.385
$ws.Range("J122").Value = 3271.5
$ws.Range("K122").Value = 18040.155
$ws.Range("L122").Value = 9814.5
$ws.Range("M122").Value = -15590.155
$ws.Range("N122").Value = -14714.5

$ws.Range("H136").Value = 1891.7391
$ws.Range("I136").Value = 1620.5
$ws.Range("K136").Value = 4861.5
$ws.Range("M136").Value = -2311.5

$ws.Range("H138").Value = 94992.5
$ws.Range("J138").Value = 94992.5
$ws.Range("L138").Value = 94992.5
$ws.Range("N138").Value = -105272.5

$ws.Range("H139").Value = 90714.664
$ws.Range("J139").Value = 90714.664
$ws.Range("L139").Value = 90714.664
$ws.Range("N139").Value = -100994.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1537402
$ws.Range("I99").Value = 60393.65
$ws.Range("J99").Value = 7814687.5
$ws.Range("K99").Value = 60393.65
$ws.Range("L99").Value = 7814687.5
$ws.Range("M99").Value = -58895.65
$ws.Range("N99").Value = -7817683.5

$ws.Range("H134").Value = 2073.64
$ws.Range("I134").Value = 1377.5
$ws.Range("K134").Value = 4132.5
$ws.Range("M134").Value = -1597.5

$ws.Range("H140").Value = 88991.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12101.385
$ws.Range("I31").Value = 2770.4707
$ws.Range("K31").Value = 2770.4707
$ws.Range("M31").Value = -2475.4707

$ws.Range("H34").Value = 12101.385
$ws.Range("I34").Value = 2770.4707
$ws.Range("K34").Value = 2770.4707
$ws.Range("M34").Value = -2568.4707

$ws.Range("H132").Value = 1641.0555
$ws.Range("I132").Value = 1518.0625
$ws.Range("J132").Value = 2625
$ws.Range("K132").Value = 4554.1875
$ws.Range("L132").Value = 7875
$ws.Range("M132").Value = -2024.1875
$ws.Range("N132").Value = -12935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2938.125
$ws.Range("J80").Value = 3166.6667
$ws.Range("L80").Value = 3166.6667
$ws.Range("N80").Value = -5162.6667

$ws.Range("H83").Value = 2938.125
$ws.Range("J83").Value = 3166.6667
$ws.Range("L83").Value = 15833.3335
$ws.Range("N83").Value = -25817.3335

$ws.Range("H102").Value = 1243.826
$ws.Range("I102").Value = 1070.5
$ws.Range("J102").Value = 2399.3333
$ws.Range("K102").Value = 1070.5
$ws.Range("L102").Value = 2399.3333
$ws.Range("M102").Value = 551.5
$ws.Range("N102").Value = -5643.3333

$ws.Range("H116").Value = 51097.777
$ws.Range("J116").Value = 51097.777
$ws.Range("L116").Value = 51097.777
$ws.Range("N116").Value = -60275.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1112.2667
$ws.Range("I16").Value = 1062.3334
$ws.Range("J16").Value = 1312
$ws.Range("K16").Value = 1062.3334
$ws.Range("L16").Value = 1312
$ws.Range("M16").Value = -892.3334
$ws.Range("N16").Value = -1652

$ws.Range("H40").Value = 3270302.8
$ws.Range("I40").Value = 2009.7273
$ws.Range("J40").Value = 9262173
$ws.Range("K40").Value = 2009.7273
$ws.Range("L40").Value = 9262173
$ws.Range("M40").Value = -1873.7273
$ws.Range("N40").Value = -9262445

$ws.Range("H61").Value = 1009.36365
$ws.Range("I61").Value = 910.3
$ws.Range("K61").Value = 910.3
$ws.Range("M61").Value = -708.3

$ws.Range("H113").Value = 1009.36365
$ws.Range("I113").Value = 910.3
$ws.Range("K113").Value = 910.3
$ws.Range("M113").Value = 1259.7

$ws.Range("H132").Value = 11171.706
$ws.Range("I132").Value = 16193.9
$ws.Range("J132").Value = 3997.1428
$ws.Range("K132").Value = 48581.7
$ws.Range("L132").Value = 11991.4284
$ws.Range("M132").Value = -46051.7
$ws.Range("N132").Value = -17051.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 145329
$ws.Range("J46").Value = 145329
$ws.Range("L46").Value = 145329
$ws.Range("N46").Value = -145791

$ws.Range("H113").Value = 696.2632
$ws.Range("I113").Value = 516.3077
$ws.Range("J113").Value = 1086.1666
$ws.Range("K113").Value = 1548.9231
$ws.Range("L113").Value = 3258.4998
$ws.Range("M113").Value = 621.0769
$ws.Range("N113").Value = -7598.4998

$ws.Range("H126").Value = 2008.0588
$ws.Range("I126").Value = 1852.1666
$ws.Range("K126").Value = 5556.4998
$ws.Range("M126").Value = -3086.4998

$ws.Range("H134").Value = 145329
$ws.Range("J134").Value = 145329
$ws.Range("L134").Value = 435987
$ws.Range("N134").Value = -441057

$ws.Range("H136").Value = 1216.4642
$ws.Range("I136").Value = 1204.56
$ws.Range("K136").Value = 3613.68
$ws.Range("M136").Value = -1063.68

$ws.Range("H137").Value = 146985
$ws.Range("J137").Value = 146985
$ws.Range("L137").Value = 146985
$ws.Range("N137").Value = -157185
